$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-autofitting the existing rows drops their stale explicit "15.75"
# height (a leftover from an older Excel build) back to the sheet
# default, and also collapses row 13 -- a height-only placeholder row
# with no cell content -- out of the sheet entirely, same as what a
# newer Excel build produces when it re-saves this file.
$ws.Rows("1:24").EntireRow.AutoFit()

# Slightly wider ID/description columns.
$ws.Columns("A").ColumnWidth = 12.83
$ws.Columns("B").ColumnWidth = 33

# --- New course-removal test case (TC007) ---
$ws.Range("A26").Value = "TC007"
$ws.Range("B27").Value = "Verify Drop course button visible"
$ws.Range("C27").Value = "1. Aftere enrolled ina  course, in list, user can remove course"
$ws.Range("D27").Value = "remove buttom is visible to remove course"

# New selection sits on the freshly added cell.
$null = $ws.Range("D27").Select()
